# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.035.13"
$ws.Cells.Item(2, 5).Value = "  +0.65%  "
$ws.Cells.Item(3, 4).Value = "2.666.05"
$ws.Cells.Item(3, 5).Value = "  +1.96%  "
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 4).Value = "'529.79"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +3.23%  "
$ws.Cells.Item(6, 4).Value = "'155.80"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  +0.77%  "
$ws.Cells.Item(7, 5).Value = "  -0.05%  "
$ws.Cells.Item(8, 4).Value = "'0.582"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -1.30%  "
$ws.Cells.Item(9, 4).Value = "'6.53"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  -4.12%  "
$ws.Cells.Item(10, 4).Value = "'0.109"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = "  +4.75%  "
$ws.Cells.Item(11, 4).Value = "'0.353"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = "  +1.73%  "
$ws.Cells.Item(12, 5).Value = "  -0.45%  "
$ws.Cells.Item(13, 4).Value = "3.135.09"
$ws.Cells.Item(13, 5).Value = "  +1.88%  "
$ws.Cells.Item(14, 4).Value = "61.023.15"
$ws.Cells.Item(14, 5).Value = "  +0.60%  "
$ws.Cells.Item(15, 4).Value = "'22.16"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "  +2.30%  "
$ws.Cells.Item(16, 5).Value = "  +1.60%  "
$ws.Cells.Item(17, 4).Value = "2.676.10"
$ws.Cells.Item(17, 5).Value = "  +1.62%  "
$ws.Cells.Item(18, 4).Value = "'4.79"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = "  +1.07%  "
$ws.Cells.Item(19, 4).Value = "'355.26"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  -0.57%  "
$ws.Cells.Item(20, 4).Value = "'10.70"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +0.79%  "
$ws.Cells.Item(21, 5).Value = "  +2.31%  "
$ws.Cells.Item(22, 4).Value = "'0.997"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -0.11%  "
$ws.Cells.Item(23, 4).Value = "'61.83"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +1.91%  "
$ws.Cells.Item(24, 5).Value = "  +1.99%  "
$ws.Cells.Item(25, 5).Value = "  +0.86%  "
$ws.Cells.Item(26, 5).Value = "  +0.05%  "
$ws.Cells.Item(27, 5).Value = "  +1.47%  "
$ws.Cells.Item(28, 5).Value = "  -0.31%  "
$ws.Cells.Item(29, 5).Value = "  -0.02%  "
$ws.Cells.Item(30, 4).Value = "'6.20"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = "  +3.94%  "
$ws.Cells.Item(31, 4).Value = "'19.55"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +0.46%  "
$ws.Cells.Item(32, 5).Value = "  +2.79%  "
$ws.Cells.Item(33, 4).Value = "'150.12"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -1.26%  "
$ws.Cells.Item(34, 4).Value = "'4.13"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +3.10%  "
$ws.Cells.Item(35, 5).Value = "  +0.37%  "
$ws.Cells.Item(36, 4).Value = "'0.926"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +8.88%  "
$ws.Cells.Item(37, 4).Value = "'0.897"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  +1.90%  "
$ws.Cells.Item(38, 2).Value = "OKB"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(38, 4).Value = "'36.89"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  +1.62%  "
$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).Value = "'306.53"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +4.93%  "
$ws.Cells.Item(40, 5).Value = "  +0.19%  "
$ws.Cells.Item(41, 5).Value = "  +0.80%  "
$ws.Cells.Item(42, 4).Value = "'0.647"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  +3.71%  "
$ws.Cells.Item(43, 5).Value = "  +0.33%  "
$ws.Cells.Item(44, 4).Value = "'20.54"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +3.36%  "
$ws.Cells.Item(45, 4).Value = "'0.0565"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = "  +1.65%  "
$ws.Cells.Item(46, 4).Value = "'0.997"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -0.08%  "
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(47, 4).Value = "'4.98"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  +1.77%  "
$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).Value = "'0.0242"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = "  +3.36%  "
$ws.Cells.Item(49, 4).Value = "'19.32"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +6.97%  "
$ws.Cells.Item(50, 4).Value = "'10.36"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = "  +0.66%  "
$ws.Cells.Item(51, 4).Value = "2.006.07"
$ws.Cells.Item(51, 5).Value = "  +0.47%  "
